# Updates cryptos list prices/volumes (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.387.82'
$ws.Range('E2').Value = '  +4.05%  '

# Row 3
$ws.Range('D3').Value = '2.431.12'
$ws.Range('E3').Value = '  +3.08%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.99'
$ws.Range('E5').Value = '  +2.40%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.37'
$ws.Range('E6').Value = '  +3.37%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  +3.52%  '

# Row 9
$ws.Range('E9').Value = '  +4.86%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.77'
$ws.Range('E10').Value = '  +3.77%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.360'
$ws.Range('E11').Value = '  +1.61%  '

# Row 12
$ws.Range('E12').Value = '  -2.47%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.99'
$ws.Range('E13').Value = '  +4.89%  '

# Row 14
$ws.Range('D14').Value = '2.863.11'
$ws.Range('E14').Value = '  +3.06%  '

# Row 15
$ws.Range('D15').Value = '60.284.37'
$ws.Range('E15').Value = '  +3.94%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000140'
$ws.Range('E16').Value = '  +4.06%  '

# Row 17
$ws.Range('D17').Value = '2.406.48'
$ws.Range('E17').Value = '  +2.38%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.44'
$ws.Range('E18').Value = '  +5.87%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.43'
$ws.Range('E19').Value = '  +3.18%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.36'
$ws.Range('E20').Value = '  +1.32%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.85'
$ws.Range('E21').Value = '  +1.86%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.52'
$ws.Range('E23').Value = '  +4.46%  '

# Row 24
$ws.Range('E24').Value = '  +3.29%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.60'
$ws.Range('E25').Value = '  +0.97%  '

# Row 26
$ws.Range('E26').Value = '  +0.06%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.36'
$ws.Range('E27').Value = '  -0.18%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0788'
$ws.Range('E28').Value = '  +6.49%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.80'
$ws.Range('E29').Value = '  +2.03%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.33'
$ws.Range('E30').Value = '  +2.97%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.25'
$ws.Range('E31').Value = '  -0.86%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.82'
$ws.Range('E32').Value = '  +1.95%  '

# Row 33
$ws.Range('E33').Value = '  +1.29%  '

# Row 34
$ws.Range('E34').Value = '  +0.00%  '

# Row 35
$ws.Range('E35').Value = '  +5.97%  '

# Row 36
$ws.Range('E36').Value = '  -0.45%  '

# Row 38
$ws.Range('E38').Value = '  +0.18%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.82'
$ws.Range('E39').Value = '  +0.94%  '

# Row 40
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.420'
$ws.Range('E40').Value = '  +10.59%  '

# Row 41
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '321.59'
$ws.Range('E41').Value = '  +10.85%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  +1.47%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.17'
$ws.Range('E43').Value = '  -0.95%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0527'
$ws.Range('E44').Value = '  +3.48%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0963'
$ws.Range('E45').Value = '  +1.89%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.65'
$ws.Range('E46').Value = '  +2.86%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.414'
$ws.Range('E47').Value = '  +8.39%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.574'
$ws.Range('E48').Value = '  +1.41%  '

# Row 49
$ws.Range('E49').Value = '  +1.57%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.97'
$ws.Range('E50').Value = '  +2.47%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  -0.12%  '
